$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows("9").Delete()
Write-Host "done"
